$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Frasca Luca"
$ws.Range("B4").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C4").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("D4").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E4").Value = "daniel pedrotti | iMontagna"
$ws.Range("F4").Value = "Davide Raffaelli | MediaserT"
